$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.318.39"
$ws.Range("E2").Value = "  +0.76%  "

# Row 3
$ws.Range("D3").Value = "1.881.62"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
$ws.Range("E4").Value = "  -0.76%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.80%  "

# Row 6
$ws.Range("E6").Value = "  -1.24%  "

# Row 7
$ws.Range("E7").Value = "  -0.75%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.355"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.48%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.61%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0749"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.77%  "

# Row 12
$ws.Range("E12").Value = "  +0.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.83%  "

# Row 14
$ws.Range("D14").Value = "2.154.12"
$ws.Range("E14").Value = "  -0.70%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.760"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.36%  "

# Row 16
$ws.Range("E16").Value = "  +0.47%  "

# Row 17
$ws.Range("D17").Value = "1.865.07"
$ws.Range("E17").Value = "  -0.88%  "

# Row 18
$ws.Range("D18").Value = "35.337.59"
$ws.Range("E18").Value = "  +0.68%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.53%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0826"
$ws.Range("E20").Value = "  +0.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.48%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.44%  "

# Row 25
$ws.Range("E25").Value = "  -0.64%  "

# Row 26
$ws.Range("E26").Value = "  -1.93%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.54%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.25%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.38%  "

# Row 30
$ws.Range("E30").Value = "  -0.02%  "

# Row 31
$ws.Range("E31").Value = "  +0.67%  "

# Row 32
$ws.Range("E32").Value = "  +1.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("E34").Value = "  -0.76%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.16%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.852"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.59%  "

# Row 38
$ws.Range("E38").Value = "  -1.93%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0727"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.65%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.77%  "

# Row 41
$ws.Range("E41").Value = "  +3.37%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "

# Row 43
$ws.Range("E43").Value = "  -1.46%  "

# Row 44
$ws.Range("E44").Value = "  +1.65%  "

# Row 45
$ws.Range("D45").Value = "1.308.95"
$ws.Range("E45").Value = "  +1.13%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0800"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.78%  "

# Row 47
$ws.Range("E47").Value = "  -1.15%  "

# Row 48
$ws.Range("E48").Value = "  -0.33%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.04%  "

# Row 50
$ws.Range("B50").Value = "Gas"
$ws.Range("C50").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.58%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "
